$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-506). The value needs to move forward by one day:
# 45177 (2023-09-08) -> 45178 (2023-09-09).
$ws.Range("C2:C506").Value = 45178
